# Fruta / hortaliza, semanal
# Insert 3 new weekly-report rows at the top of the data block (row 9),
# pushing the existing rows down, then populate the new rows with the
# latest week's figures (Region de O'Higgins, week of 2023-06-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 9:32 down to 12:35.
$ws.Rows("9:11").Insert()

function Set-Row {
    param($r, $d, $l, $m, $n, $o, $p, $rOrigen, $s)

    $ws.Cells.Item($r, 1).Value = 3
    $ws.Cells.Item($r, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = 5
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100107
    $ws.Cells.Item($r, 8).Value = "Otros"
    $ws.Cells.Item($r, 9).Value = 100107001
    $ws.Cells.Item($r, 10).Value = "Caqui"
    $ws.Cells.Item($r, 11).Value = "Mankaki"
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 10 kilos"
    $ws.Cells.Item($r, 18).Value = $rOrigen
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = 10
}

Set-Row 9  45082 "Especial" 56 15000 15000 15000 "Región de O'Higgins" 1500
Set-Row 10 45082 "Primera"  67 12000 12000 12000 "Región de O'Higgins" 1200
Set-Row 11 45082 "Segunda"  60 10000 10000 10000 "Región de O'Higgins" 1000

Write-Output "done"
